$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 33, shifting existing rows 33-60 down to 35-62.
$ws.Rows.Item(33).Resize(2).Insert()

# Row 33 (new): Primera, 27-dec-2021
$ws.Cells.Item(33, 1).Value = 9
$ws.Cells.Item(33, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(33, 3).Value = "Metropolitana"
$ws.Cells.Item(33, 4).Value = 44557
$ws.Cells.Item(33, 5).Value = 13
$ws.Cells.Item(33, 6).Value = 100114002
$ws.Cells.Item(33, 7).Value = "Camote"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 1600
$ws.Cells.Item(33, 11).Value = 11000
$ws.Cells.Item(33, 12).Value = 12000
$ws.Cells.Item(33, 13).Value = 11500
$ws.Cells.Item(33, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(33, 15).Value = "Perú"
$ws.Cells.Item(33, 16).Value = 639
$ws.Cells.Item(33, 17).Value = 18
$ws.Cells.Item(33, 18).Value = "Hortaliza"

# Row 34 (new): Segunda
$ws.Cells.Item(34, 1).Value = 9
$ws.Cells.Item(34, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(34, 3).Value = "Metropolitana"
$ws.Cells.Item(34, 4).Value = 44557
$ws.Cells.Item(34, 5).Value = 13
$ws.Cells.Item(34, 6).Value = 100114002
$ws.Cells.Item(34, 7).Value = "Camote"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Segunda"
$ws.Cells.Item(34, 10).Value = 790
$ws.Cells.Item(34, 11).Value = 9000
$ws.Cells.Item(34, 12).Value = 9000
$ws.Cells.Item(34, 13).Value = 9000
$ws.Cells.Item(34, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(34, 15).Value = "Perú"
$ws.Cells.Item(34, 16).Value = 500
$ws.Cells.Item(34, 17).Value = 18
$ws.Cells.Item(34, 18).Value = "Hortaliza"

# Apply the same date-number style as the other Fecha (column D) cells.
$ws.Cells.Item(33, 4).NumberFormat = $ws.Cells.Item(35, 4).NumberFormat
$ws.Cells.Item(34, 4).NumberFormat = $ws.Cells.Item(35, 4).NumberFormat
